# Regenerate save_data: use K (strikeouts, recalculated) instead of the
# previous Strike# values in column G. Updates rows 2-28 of the "K" column.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$newValues = @{
    2  = 5
    3  = 5
    4  = 4
    5  = 5
    6  = 2
    7  = 7
    8  = 6
    9  = 5
    10 = 5
    11 = 2
    12 = 10
    13 = 5
    14 = 4
    15 = 5
    16 = 1
    17 = 1
    18 = 4
    19 = 2
    20 = 3
    21 = 5
    22 = 4
    23 = 2
    24 = 5
    25 = 3
    26 = 3
    27 = 3
    28 = 1
}

foreach ($row in $newValues.Keys) {
    $ws.Range("G$row").Value = $newValues[$row]
}
